$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.339.41'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.57'
$ws.Range('E3').Value = '  +2.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.81'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4450'
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3767'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.66'
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07740'
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.133'
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.30'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.339'
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.570'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.843.41'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.50'
$ws.Range('E17').Value = '  +15.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001084'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06492'
$ws.Range('E19').Value = '  -3.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.55'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.349'
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5421'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '28.385.19'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.73'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.209'
$ws.Range('E26').Value = '  -9.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.77'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.61'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.366'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.043.87'
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.67'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.205'
$ws.Range('E32').Value = '  -10.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.921'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09299'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.673'
$ws.Range('E35').Value = '  -7.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '13.20'
$ws.Range('E36').Value = '  +8.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02355'
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2195'
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.207'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6612'
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06235'
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.203'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.153'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.98'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.391'
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6117'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.783'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.056'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.51'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.156'
$ws.Range('E51').Value = '  +1.29%  '
